# Test with a variable production cost
# Update the StartingInventories (column C) values on the "Productdata" sheet
# for Part_0001, Part_0002, Part_0003, Retail_0001, Retail_0002, Retail_0003
# (rows 4-9) from 0 to their new non-zero values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Productdata")

$ws.Range("C4").Value = 1673
$ws.Range("C5").Value = 1673
$ws.Range("C6").Value = 1673
$ws.Range("C7").Value = 1013
$ws.Range("C8").Value = 181
$ws.Range("C9").Value = 301
